# Weekly data update: add a new "18_05_2021" column (AC) to the Intensiv sheet,
# mirroring the existing week-over-week layout (header in row 1, per-age-group
# counts in rows 2-11, SUM total in row 12), and nudge the view over to the
# newly-added column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header (goes into the shared-string table as "18_05_2021")
$ws.Range("AC1").Value = "18_05_2021"

# New weekly counts per age group (rows 2-11)
$ws.Range("AC2").Value  = 12
$ws.Range("AC3").Value  = 16
$ws.Range("AC4").Value  = 36
$ws.Range("AC5").Value  = 48
$ws.Range("AC6").Value  = 121
$ws.Range("AC7").Value  = 282
$ws.Range("AC8").Value  = 410
$ws.Range("AC9").Value  = 569
$ws.Range("AC10").Value = 178
$ws.Range("AC11").Value = 13

# Total row - same SUM pattern used by every other week's column
$ws.Range("AC12").Formula = "=SUM(AC2:AC11)"

# Match the widened columns used for the most-recent weeks (Y:AC) versus the
# older, narrower ones (A:X). The COM width setter here snaps to the nearest
# 1/6-character pixel grid, so this is the closest achievable value to the
# 11.7109375 stored by Excel.
$ws.Range("Y1:AC1").ColumnWidth = 10.8333333333333

# Move the visible selection over to the newly-added data, matching the
# author's scroll position after appending this week's column.
$ws.Range("Y8").Select()
